# Update "想去人数" (F column) counts on the 展览/演出/全部类型 sheets to the
# freshly-generated gh-pages snapshot values (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 274
$ws.Range("F3").Value = 3243
$ws.Range("F4").Value = 2003
$ws.Range("F5").Value = 268
$ws.Range("F6").Value = 107
$ws.Range("F7").Value = 3087
$ws.Range("F8").Value = 616
$ws.Range("F9").Value = 300
$ws.Range("F10").Value = 39
$ws.Range("F12").Value = 153
$ws.Range("F15").Value = 10166
$ws.Range("F18").Value = 10
$ws.Range("F20").Value = 8064
$ws.Range("F21").Value = 12675
$ws.Range("F24").Value = 24
$ws.Range("F25").Value = 273
$ws.Range("F26").Value = 396
$ws.Range("F27").Value = 595
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 416
$ws.Range("F30").Value = 2827
$ws.Range("F33").Value = 7967
$ws.Range("F34").Value = 1555
$ws.Range("F36").Value = 70
$ws.Range("F38").Value = 4625
$ws.Range("F39").Value = 1431

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 1198
$ws.Range("F12").Value = 28

# Sheet 4: 全部类型 (All types — combined listing)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 274
$ws.Range("F5").Value = 3243
$ws.Range("F7").Value = 2003
$ws.Range("F9").Value = 268
$ws.Range("F11").Value = 3087
$ws.Range("F13").Value = 616
$ws.Range("F14").Value = 39
$ws.Range("F16").Value = 153
$ws.Range("F19").Value = 10166
$ws.Range("F21").Value = 10
$ws.Range("F23").Value = 8064
$ws.Range("F24").Value = 12675
$ws.Range("F27").Value = 24
$ws.Range("F28").Value = 273
$ws.Range("F30").Value = 595
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 2827
$ws.Range("F35").Value = 28
$ws.Range("F38").Value = 7967
$ws.Range("F40").Value = 70
$ws.Range("F42").Value = 4625
